$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.75"
$ws.Range("E2").Value = "'-4.01%"
$ws.Range("D3").Value = "'37.09"
$ws.Range("E3").Value = "'-6.99%"
$ws.Range("D4").Value = "'5.084"
$ws.Range("E4").Value = "'-1.15%"
$ws.Range("D5").Value = "'0.07708"
$ws.Range("E5").Value = "'-6.18%"
$ws.Range("D6").Value = "'4.349"
$ws.Range("E6").Value = "'0.68%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.170"
$ws.Range("E7").Value = "'-2.58%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.877"
$ws.Range("E8").Value = "'-8.91%"
$ws.Range("E9").Value = "'-4.90%"
$ws.Range("D10").Value = "'0.9209"
$ws.Range("E10").Value = "'-2.23%"
$ws.Range("D11").Value = "'0.1200"
$ws.Range("E11").Value = "'-11.30%"
$ws.Range("D12").Value = "'0.1856"
$ws.Range("E12").Value = "'-7.11%"
$ws.Range("D13").Value = "'0.08754"
$ws.Range("E13").Value = "'-5.32%"
$ws.Range("D14").Value = "'0.03391"
$ws.Range("E14").Value = "'-3.78%"
$ws.Range("D15").Value = "'0.09683"
$ws.Range("E15").Value = "'-1.12%"
$ws.Range("D16").Value = "'0.001377"
$ws.Range("E16").Value = "'-2.35%"
$ws.Range("D17").Value = "'0.006004"
$ws.Range("E17").Value = "'-6.55%"
$ws.Range("D18").Value = "'3.603"
$ws.Range("E18").Value = "'-2.31%"
$ws.Range("E19").Value = "'-2.42%"
$ws.Range("E20").Value = "'-3.44%"
$ws.Range("D21").Value = "'5.015"
$ws.Range("E21").Value = "'1.14%"
$ws.Range("E22").Value = "'6.04%"
$ws.Range("D23").Value = "'0.02107"
$ws.Range("E23").Value = "'5,163.64%"
$ws.Range("D24").Value = "'0.04324"
$ws.Range("E24").Value = "'-0.90%"
$ws.Range("D25").Value = "'0.001213"
$ws.Range("E25").Value = "'-1.64%"
$ws.Range("D26").Value = "'0.004209"
$ws.Range("E26").Value = "'-12.11%"
$ws.Range("D27").Value = "'0.0001352"
$ws.Range("E27").Value = "'3.86%"
$ws.Range("D39").Value = "'0.02172"
$ws.Range("E39").Value = "'-6.86%"
$ws.Range("D40").Value = "'0.04883"
$ws.Range("E40").Value = "'-6.29%"
$ws.Range("D41").Value = "'0.007568"
$ws.Range("E41").Value = "'-2.45%"
$ws.Range("D42").Value = "'0.009938"
$ws.Range("E42").Value = "'0.39%"
$ws.Range("D43").Value = "'0.1337"
$ws.Range("E43").Value = "'-4.94%"
$ws.Range("D44").Value = "'0.001996"
$ws.Range("E44").Value = "'-4.11%"
$ws.Range("D45").Value = "'0.009112"
$ws.Range("E45").Value = "'-1.09%"
$ws.Range("D46").Value = "'0.00006513"
$ws.Range("E46").Value = "'-1.60%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.09%"
$ws.Range("D48").Value = "'0.003003"
$ws.Range("E48").Value = "'1.91%"
$ws.Range("D49").Value = "'0.001303"
$ws.Range("E49").Value = "'-22.97%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.09%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.09%"
